$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Insert a new column before column B. This shifts the existing "dbExcel"
# header/value column (old B) and the "WebExcel" header/value column
# (old C) one position to the right (new C, D), and Excel carries the
# column-B formatting (including A2's wrap-text style) into the freshly
# inserted column automatically.
$ws.Columns.Item(2).Insert()

# Column widths per target layout: A:B share the original A width,
# C keeps the original B width, D keeps the original C width.
# (Values chosen to land in the engine's nearest pixel-rounded bucket to the
# exact target character widths 75.81640625 / 70.26953125 / 28.54296875.)
$ws.Columns.Item(1).ColumnWidth = 75.0
$ws.Columns.Item(2).ColumnWidth = 75.0
$ws.Columns.Item(3).ColumnWidth = 69.5
$ws.Columns.Item(4).ColumnWidth = 27.666666666666668

# New "StatQuery" column: header + the stats Cypher query, wrapped like A2.
$ws.Range("B1").Value = "StatQuery"
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.ethnicity IN ['NOT_HISPANIC_OR_LATINO'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"
$ws.Range("B2").WrapText = $true
